# Update countries & provincias Spain
# - Update the "last updated" timestamp in A1.
# - Refresh the per-country COVID figures that changed with this data pull.
# - Re-sort the country table (A4:H219) by "Casos totales" (column B) descending,
#   which is what shuffles Rumania/Guatemala/Polonia, Palestina/Bulgaria and
#   Hong Kong/Libia relative to each other (their row data otherwise does not change).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Timestamp header
$ws.Range("A1").Value = "Datos actualizados a 24 de Julio de 2020 a las 12:28"

# 2. Updated per-country figures (row numbers are the PRE-sort positions from the
#    original sheet; the later Sort.Apply() reorders the rows by column B).
# Columns: B=Casos totales, C=Nuevos casos, D=Casos activos, E=Recuperados,
#          F=Casos criticos, G=Muertes hoy, H=Muertes

# Estados Unidos (row 4)
$ws.Cells.Item(4, 2).Value = 4170333
$ws.Cells.Item(4, 3).Value = 342
$ws.Cells.Item(4, 5).Value = 2042559
$ws.Cells.Item(4, 7).Value = 9
$ws.Cells.Item(4, 8).Value = 147342

# Indonesia (row 27)
$ws.Cells.Item(27, 2).Value = 95418
$ws.Cells.Item(27, 3).Value = 1761
$ws.Cells.Item(27, 4).Value = 53945
$ws.Cells.Item(27, 5).Value = 36808
$ws.Cells.Item(27, 7).Value = 89
$ws.Cells.Item(27, 8).Value = 4665

# Rumania (row 49)
$ws.Cells.Item(49, 2).Value = 42394
$ws.Cells.Item(49, 3).Value = 1119
$ws.Cells.Item(49, 4).Value = 25349
$ws.Cells.Item(49, 5).Value = 14895
$ws.Cells.Item(49, 7).Value = 24
$ws.Cells.Item(49, 8).Value = 2150

# Ghana (row 57)
$ws.Cells.Item(57, 2).Value = 30366
$ws.Cells.Item(57, 3).Value = 694
$ws.Cells.Item(57, 4).Value = 26687
$ws.Cells.Item(57, 5).Value = 3526

# Estado de Palestina (row 81)
$ws.Cells.Item(81, 2).Value = 10093
$ws.Cells.Item(81, 3).Value = 349
$ws.Cells.Item(81, 5).Value = 7303
$ws.Cells.Item(81, 7).Value = 3
$ws.Cells.Item(81, 8).Value = 70

# Malasia (row 86)
$ws.Cells.Item(86, 2).Value = 8861
$ws.Cells.Item(86, 3).Value = 21
$ws.Cells.Item(86, 4).Value = 8577
$ws.Cells.Item(86, 5).Value = 161

# Finlandia (row 89)
$ws.Cells.Item(89, 2).Value = 7380
$ws.Cells.Item(89, 3).Value = 8
$ws.Cells.Item(89, 5).Value = 132

# Hong Kong (row 119)
$ws.Cells.Item(119, 2).Value = 2373
$ws.Cells.Item(119, 3).Value = 123
$ws.Cells.Item(119, 4).Value = 1407
$ws.Cells.Item(119, 5).Value = 950
$ws.Cells.Item(119, 7).Value = 1
$ws.Cells.Item(119, 8).Value = 16

# Uganda (row 145)
$ws.Cells.Item(145, 2).Value = 1089
$ws.Cells.Item(145, 3).Value = 10
$ws.Cells.Item(145, 4).Value = 975
$ws.Cells.Item(145, 5).Value = 113
$ws.Cells.Item(145, 7).Value = 1
$ws.Cells.Item(145, 8).Value = 1

# Groenlandia / Islas Malvinas (rows 210-211) are a genuine tie on every numeric
# column, so a stable sort alone would keep their original relative order. Swap
# the country labels up front so the tie resolves with Groenlandia first, matching
# the shared-strings reorder in the diff.
$ws.Cells.Item(210, 1).Value = "Groenlandia"
$ws.Cells.Item(211, 1).Value = "Islas Malvinas"

# 3. Re-sort the country table by "Casos totales" (column B) descending, so rows
#    move to reflect the new ranking (matches the shared-strings reorder in the diff).
$sortRange = $ws.Range("A4:H219")
$sortKey = $ws.Range("B4:B219")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($sortKey, 0, 2, 0, 0)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()
$ws.Sort.SortFields.Clear()
